$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "319.85"
    "E2" = "3.48%"
    "D3" = "41.39"
    "E3" = "1.37%"
    "D4" = "5.260"
    "E4" = "2.77%"
    "B6" = "GateToken"
    "C6" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D6" = "4.341"
    "E6" = "1.46%"
    "B7" = "FTXToken"
    "C7" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D7" = "1.759"
    "E7" = "9.76%"
    "B8" = "MXToken"
    "C8" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8" = "0.9437"
    "E8" = "3.96%"
    "B9" = "BTSEToken"
    "C9" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D9" = "2.425"
    "E9" = "-1.58%"
    "B10" = "LiechtensteinCryptoassetsExchange"
    "C10" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D10" = "0.1244"
    "E10" = "-2.32%"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D11" = "0.1878"
    "E11" = "4.06%"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D12" = "0.09179"
    "E12" = "1.50%"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D13" = "0.04308"
    "E13" = "0.35%"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "0.1052"
    "E14" = "0.67%"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "0.001296"
    "E15" = "3.72%"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D16" = "0.005869"
    "E16" = "3.95%"
    "B17" = "UpBots"
    "C17" = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
    "D17" = "0.007491"
    "E17" = "1,897.31%"
    "B18" = "LEO"
    "C18" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D18" = "3.342"
    "E18" = "-0.22%"
    "D19" = "0.3357"
    "E19" = "1.30%"
    "D20" = "7.790"
    "E20" = "12.68%"
    "D21" = "0.1344"
    "E21" = "-3.57%"
    "D22" = "0.2824"
    "E22" = "4.37%"
    "D23" = "0.04040"
    "E23" = "0.00%"
    "E24" = "-0.23%"
    "D25" = "0.004120"
    "E25" = "1.96%"
    "D26" = "0.0001272"
    "E26" = "-0.01%"
    "D38" = "0.02554"
    "E38" = "5.70%"
    "D39" = "0.05348"
    "E39" = "2.40%"
    "D40" = "0.007766"
    "E40" = "-1.03%"
    "D41" = "0.1319"
    "E41" = "1.49%"
    "D42" = "0.007051"
    "E42" = "3.60%"
    "D43" = "0.001992"
    "E43" = "3.00%"
    "D44" = "0.008275"
    "E44" = "12.70%"
    "E45" = "-5.29%"
    "D46" = "0.00006700"
    "E46" = "-2.80%"
    "E47" = "-0.02%"
    "D48" = "0.1994"
    "E48" = "51.62%"
    "D49" = "0.004207"
    "E49" = "40.12%"
    "E50" = "-0.02%"
    "D51" = "0.0002003"
    "E51" = "-0.02%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
